$d = $word.ActiveDocument

# 1. Remove the leading curly opening quote (U+201C) from the start of the
#    "Generally, symbols are placed..." paragraph's text. The rest of the
#    text (including the later quoted "IM") is left untouched.
$openQuote = [char]0x201C
$d.Content.Find.Execute($openQuote + "Generally, symbols", $true, $false, $false, $false, $false, $true, 1, $false, "Generally, symbols", 2) | Out-Null

# 2. Move the "_GoBack" bookmark from the very start of the document (before
#    the "Kata: Roman Numerals" run) to the very start of the paragraph that
#    now begins "Generally, symbols...".
$bookmarkName = "_GoBack"
if ($d.Bookmarks.Exists($bookmarkName)) {
    $d.Bookmarks.Item($bookmarkName).Delete()
}

$findRange = $d.Content
$findRange.Find.Execute("Generally, symbols are placed") | Out-Null
$targetRange = $d.Range($findRange.Start, $findRange.Start)
$d.Bookmarks.Add($bookmarkName, $targetRange) | Out-Null
